# Update "想去人数" (F column) counts across the four sheets to the newly
# scraped values, as recorded by the gh-pages output regeneration commit.

$wb = $excel.ActiveWorkbook

# sheet name -> hashtable of Row -> NewValue
$changesBySheet = @{
    "展览" = @{
        4  = 638
        5  = 479
        6  = 262
        7  = 1148
        9  = 177
        10 = 65
        11 = 764
        12 = 410
        13 = 55
        15 = 201
        18 = 6241
        19 = 51
        22 = 7214
        25 = 3292
        26 = 418
        27 = 794
        28 = 4476
        29 = 332
        30 = 157
        31 = 157
        32 = 1281
        33 = 119
        34 = 38
        36 = 999
        37 = 1301
        38 = 2074
    }
    "演出" = @{
        2 = 56
    }
    "本地生活" = @{
        3 = 1169
        4 = 60
    }
    "全部类型" = @{
        4  = 1169
        5  = 60
        7  = 638
        8  = 479
        9  = 262
        10 = 1148
        12 = 177
        13 = 65
        14 = 764
        15 = 410
        16 = 55
        17 = 56
        19 = 202
        22 = 6241
        23 = 6241
        24 = 51
        27 = 7214
        30 = 3292
        31 = 418
        32 = 794
        33 = 4476
        34 = 332
        36 = 157
        37 = 157
        38 = 1281
        39 = 119
        40 = 38
        42 = 999
        43 = 1301
        45 = 2074
    }
}

foreach ($sheetName in $changesBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $changesBySheet[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}
